$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates reflecting the refreshed crypto price/volume snapshot.
# NumberFormat is forced to Text ("@") before each write so that numeric-
# looking strings (e.g. "1.002", "0.000008940") are preserved verbatim
# instead of being auto-coerced into floating point numbers.
$updates = @(
    @{Cell='D2'; Value='27.656.46'}
    @{Cell='E2'; Value='  +0.24%  '}
    @{Cell='D3'; Value='1.845.27'}
    @{Cell='E3'; Value='  +0.22%  '}
    @{Cell='D4'; Value='1.002'}
    @{Cell='E4'; Value='  +0.11%  '}
    @{Cell='D5'; Value='312.61'}
    @{Cell='E5'; Value='  -0.48%  '}
    @{Cell='E6'; Value='  +0.01%  '}
    @{Cell='D7'; Value='0.4277'}
    @{Cell='E7'; Value='  +0.89%  '}
    @{Cell='D8'; Value='0.3631'}
    @{Cell='E8'; Value='  -0.06%  '}
    @{Cell='D9'; Value='0.07328'}
    @{Cell='E9'; Value='  +0.96%  '}
    @{Cell='D10'; Value='0.8766'}
    @{Cell='D11'; Value='20.62'}
    @{Cell='E11'; Value='  +0.30%  '}
    @{Cell='D12'; Value='1.871.96'}
    @{Cell='E12'; Value='  +1.57%  '}
    @{Cell='D13'; Value='5.353'}
    @{Cell='E13'; Value='  -0.05%  '}
    @{Cell='D14'; Value='6.519'}
    @{Cell='E14'; Value='  -0.55%  '}
    @{Cell='D15'; Value='0.06952'}
    @{Cell='E15'; Value='  +1.35%  '}
    @{Cell='D16'; Value='1.004'}
    @{Cell='E16'; Value='  +0.15%  '}
    @{Cell='D17'; Value='79.57'}
    @{Cell='E17'; Value='  +1.45%  '}
    @{Cell='D18'; Value='0.000008940'}
    @{Cell='E18'; Value='  +1.64%  '}
    @{Cell='D19'; Value='1.003'}
    @{Cell='E19'; Value='  +0.17%  '}
    @{Cell='E20'; Value='  -0.51%  '}
    @{Cell='D21'; Value='27.807.74'}
    @{Cell='E21'; Value='  +0.86%  '}
    @{Cell='D22'; Value='4.985'}
    @{Cell='E22'; Value='  -0.10%  '}
    @{Cell='E23'; Value='  -2.21%  '}
    @{Cell='D24'; Value='2.131.07'}
    @{Cell='E24'; Value='  +3.14%  '}
    @{Cell='E25'; Value='  -2.13%  '}
    @{Cell='D26'; Value='155.48'}
    @{Cell='E26'; Value='  +0.21%  '}
    @{Cell='D27'; Value='18.50'}
    @{Cell='E27'; Value='  -0.16%  '}
    @{Cell='D28'; Value='119.60'}
    @{Cell='E28'; Value='  +0.34%  '}
    @{Cell='D29'; Value='5.208'}
    @{Cell='E29'; Value='  +0.04%  '}
    @{Cell='D30'; Value='1.869'}
    @{Cell='E30'; Value='  +2.88%  '}
    @{Cell='D31'; Value='0.08873'}
    @{Cell='E31'; Value='  +0.01%  '}
    @{Cell='D32'; Value='0.7544'}
    @{Cell='E32'; Value='  -2.43%  '}
    @{Cell='D33'; Value='2.964'}
    @{Cell='E33'; Value='  +0.41%  '}
    @{Cell='E34'; Value='  -0.88%  '}
    @{Cell='D35'; Value='1.128'}
    @{Cell='E35'; Value='  +2.73%  '}
    @{Cell='E36'; Value='  +0.04%  '}
    @{Cell='E37'; Value='  +0.73%  '}
    @{Cell='D38'; Value='1.107'}
    @{Cell='E38'; Value='  +0.98%  '}
    @{Cell='D39'; Value='0.01932'}
    @{Cell='E39'; Value='  +0.68%  '}
    @{Cell='D40'; Value='2.820'}
    @{Cell='E40'; Value='  +1.93%  '}
    @{Cell='D41'; Value='0.1664'}
    @{Cell='E41'; Value='  +1.09%  '}
    @{Cell='D42'; Value='0.5069'}
    @{Cell='E42'; Value='  +0.42%  '}
    @{Cell='D43'; Value='6.591'}
    @{Cell='E43'; Value='  -3.19%  '}
    @{Cell='D44'; Value='8.387'}
    @{Cell='E44'; Value='  +2.57%  '}
    @{Cell='D45'; Value='0.06546'}
    @{Cell='E45'; Value='  -0.94%  '}
    @{Cell='B46'; Value='Quant'}
    @{Cell='C46'; Value='https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'}
    @{Cell='D46'; Value='106.15'}
    @{Cell='E46'; Value='  +1.00%  '}
    @{Cell='B47'; Value='EnergySwap'}
    @{Cell='C47'; Value='https://coinranking.com/coin/SbWqqTui-+energyswap-ens'}
    @{Cell='D47'; Value='10.36'}
    @{Cell='E47'; Value='  +0.91%  '}
    @{Cell='E48'; Value='  -0.83%  '}
    @{Cell='D49'; Value='1.001'}
    @{Cell='E49'; Value='  +0.05%  '}
    @{Cell='D50'; Value='1.637'}
    @{Cell='E50'; Value='  +1.06%  '}
    @{Cell='D51'; Value='64.74'}
    @{Cell='E51'; Value='  +0.50%  '}
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
}
